$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Day of Judgment', ['{2}{W}{W}', 'Sorcery', 'Destroy all creatures.'])"
$ws.Range("A3").Value = "('Emeria Angel', ['{2}{W}{W}', 'Creature — Angel', 'Flying', 'Landfall — Whenever a land enters the battlefield under your control, you may create a 1/1 white Bird creature token with flying.', '3/3'])"
$ws.Range("A4").Value = "(`"Nissa's Chosen`", ['{G}{G}', 'Creature — Elf Warrior', 'If Nissa’s Chosen would die, put it on the bottom of its owner’s library instead.', '2/3'])"
$ws.Range("A5").Value = "('Rampaging Baloths', ['{4}{G}{G}', 'Creature — Beast', 'Trample', 'Landfall — Whenever a land enters the battlefield under your control, you may create a 4/4 green Beast creature token.', '6/6'])"
$ws.Range("A6").Value = "('Valakut, the Molten Pinnacle', ['Land', 'Valakut, the Molten Pinnacle enters the battlefield tapped.', 'Whenever a Mountain enters the battlefield under your control, if you control at least five other Mountains, you may have Valakut, the Molten Pinnacle deal 3 damage to any target.', '{T}: Add {R}.'])"

$ws.Range("A7:A27").EntireRow.Delete()
